# Auto update: 2025-12-05 18:24:43
# Update the K (최종점수) and N (MACRO_SCORE) columns with refreshed scores,
# and correct the M2 (판단) decision label for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 최종점수 (final score) column K, rows 2-5
$ws.Range("K2").Value = 59.9
$ws.Range("K3").Value = 57.5
$ws.Range("K4").Value = 50.7
$ws.Range("K5").Value = 47.9

# 판단 (decision) for UnitedHealth (row 2) moved from "buy watch" to "stand aside"
$ws.Range("M2").Value = "⛔ 관망하십시오."

# MACRO_SCORE column N, rows 2-5, refreshed value
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("N3").Value = 51.53902399942638
$ws.Range("N4").Value = 51.53902399942638
$ws.Range("N5").Value = 51.53902399942638
